$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting from column N (previous year) into the new column O
$ws.Range("N4").Copy()
$ws.Range("O4").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("N5").Copy()
$ws.Range("O5").PasteSpecial(-4122)  # xlPasteFormats

# Add new column O data for year 2021
$ws.Range("O4").Value = 2021
$ws.Range("O5").Value = 515

# Update the active cell selection view
$ws.Range("P12").Select()
